# Actualización automática 2025-07-02 14:35:07
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": LAVABOS figures for SOLIS OCAMPO DIMAS ABDON ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("I9").Value = 43.86
$wsGrupo.Range("I10").Value = "1 de 8"

# --- Sheet "VENTA MENSUAL": julio figures for SOLIS OCAMPO DIMAS ABDON ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F9").Value = 43.86
$wsMensual.Range("F10").Value = 110.06
